$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the two threshold values that changed
$ws.Range("B2").Value = 5.5
$ws.Range("B3").Value = 6

# Move/restore the active selection to D6 (single cell)
$null = $ws.Range("D6").Select()
